$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row of data at row 15 (row 14 stays empty)
$ws.Range("A15").Value = 20250721
$ws.Range("B15").Value = 3
$ws.Range("C15").Value = 6
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 5

# Update the selected cell shown in the sheet view
$ws.Range("D18").Select()
